$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing total_venda (B column) values for 07/2025 entries
$ws.Range("B2").Value = 18041.08
$ws.Range("B12").Value = 48059.45
$ws.Range("B13").Value = 16959.9

# Insert a new row for dia 17 (07/2025) right after the existing dia 16 row,
# shifting the remaining rows (dia 2..30 of 06/2025, 05/2025, 04/2025) down by one
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 17
$ws.Range("B14").Value = 12257.2
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 2025
$ws.Range("E14").Value = "07/2025"
